# Update the "想去人数" (F column) counts that changed between data refreshes.
# Sheet "展览" (Exhibitions)
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 8020
$ws1.Range("F5").Value  = 5849
$ws1.Range("F6").Value  = 497
$ws1.Range("F7").Value  = 85
$ws1.Range("F10").Value = 285
$ws1.Range("F11").Value = 377

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 90

# Sheet "全部类型" (All types - aggregated view)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 8020
$ws4.Range("F5").Value  = 5849
$ws4.Range("F6").Value  = 497
$ws4.Range("F7").Value  = 85
$ws4.Range("F10").Value = 285
$ws4.Range("F11").Value = 90
$ws4.Range("F14").Value = 377
